# The "2024" worksheet holds a running log of bank/SMS notifications grouped
# by category in columns A (group label) plus paired Details/Date columns
# per month (September lives in columns R/S). A brand-new September entry
# ("dispute" at 2024-09-09 12:17:30) was logged, which pushes every row from
# 35 downward (through the August R/S->P/Q tail and the trailing "Others"/
# "Broadband" group markers) down by one row, growing the sheet from
# A1:Y116 to A1:Y117.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a fresh blank row at 35 - this shifts rows 35:116 down to 36:117,
# which reproduces every cascading R/S (and eventually P/Q, A) shift seen
# in the diff without touching any of the other, unrelated rows.
$ws.Rows(35).Insert()

# Populate the newly-opened row 35 with the new September notification.
$ws.Range("R35").Value = "dispute"
$ws.Range("S35").Value = "2024-09-09 12:17:30"
